# SystemTest/SmokeSanityTestCases.xlsx update for LPH_V1.9.2
# - Adds "Deprecated: LPH" / "Hidden: LPH" / "OnHold: LPH" labels in column F
#   for rows 14-18 (TC016-TC019 area of the sanity test cases).
# - Moves the active selection on the sheet from A23 to D21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F14").Value = "Deprecated: LPH"
$ws.Range("F15").Value = "Hidden: LPH"
$ws.Range("F16").Value = "OnHold: LPH"
$ws.Range("F17").Value = "OnHold: LPH"
$ws.Range("F18").Value = "OnHold: LPH"

# Update the saved selection/active cell for the sheet.
$ws.Range("D21").Select()
